# Auto-generated PowerShell-style Excel COM-interop edit script
# EIA Table 2.13.B 2017-01-31 update: October 2016/2015 YTD vintage -> November 2016/2015 YTD vintage

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title (A2) text update: October -> November ---
$ws.Range("A2").Value = "Year-to-Date through November 2016 and November 2015 (Thousand Tons)"

# --- Column-header (row 5) text updates: "October 2016 YTD" -> "November 2016 YTD" and "October 2015 YTD" -> "November 2015 YTD" ---
$ws.Range("B5").Value = "November 2016 YTD"
$ws.Range("E5").Value = "November 2016 YTD"
$ws.Range("G5").Value = "November 2016 YTD"
$ws.Range("I5").Value = "November 2016 YTD"
$ws.Range("K5").Value = "November 2016 YTD"
$ws.Range("C5").Value = "November 2015 YTD"
$ws.Range("F5").Value = "November 2015 YTD"
$ws.Range("H5").Value = "November 2015 YTD"
$ws.Range("J5").Value = "November 2015 YTD"
$ws.Range("L5").Value = "November 2015 YTD"

# --- State / region data updates (values refreshed to the November 2016/2015 YTD vintage) ---
# Row 6
$ws.Range("B6").Value = 3253
$ws.Range("C6").Value = 3440
$ws.Range("D6").Value = -0.054
$ws.Range("G6").Value = 3084
$ws.Range("H6").Value = 3246
$ws.Range("I6").Value = 170
$ws.Range("J6").Value = 195

# Row 7
$ws.Range("B7").Value = 1146
$ws.Range("C7").Value = 1180
$ws.Range("D7").Value = -0.029
$ws.Range("G7").Value = 1146
$ws.Range("H7").Value = 1165

# Row 8
$ws.Range("B8").Value = 251
$ws.Range("C8").Value = 271
$ws.Range("D8").Value = -0.074
$ws.Range("G8").Value = 81
$ws.Range("H8").Value = 91
$ws.Range("I8").Value = 170
$ws.Range("J8").Value = 180

# Row 9
$ws.Range("B9").Value = 1759
$ws.Range("C9").Value = 1877
$ws.Range("D9").Value = -0.063
$ws.Range("G9").Value = 1759
$ws.Range("H9").Value = 1877

# Row 10
$ws.Range("B10").Value = 97
$ws.Range("C10").Value = 112
$ws.Range("G10").Value = 97
$ws.Range("H10").Value = 112

# Row 13
$ws.Range("B13").Value = 4613
$ws.Range("C13").Value = 4860
$ws.Range("D13").Value = -0.051
$ws.Range("G13").Value = 3691
$ws.Range("H13").Value = 3835
$ws.Range("I13").Value = 922
$ws.Range("J13").Value = 1025

# Row 14
$ws.Range("B14").Value = 1219
$ws.Range("C14").Value = 1253
$ws.Range("D14").Value = -0.027
$ws.Range("G14").Value = 905
$ws.Range("H14").Value = 937
$ws.Range("I14").Value = 314
$ws.Range("J14").Value = 316

# Row 15
$ws.Range("B15").Value = 1696
$ws.Range("C15").Value = 1839
$ws.Range("D15").Value = -0.078
$ws.Range("G15").Value = 1311
$ws.Range("H15").Value = 1373
$ws.Range("I15").Value = 385
$ws.Range("J15").Value = 466

# Row 16
$ws.Range("B16").Value = 1698
$ws.Range("C16").Value = 1767
$ws.Range("D16").Value = -0.039
$ws.Range("G16").Value = 1475
$ws.Range("H16").Value = 1525
$ws.Range("I16").Value = 223
$ws.Range("J16").Value = 243

# Row 17
$ws.Range("B17").Value = 226
$ws.Range("C17").Value = 237
$ws.Range("D17").Value = -0.044
$ws.Range("E17").Value = 36
$ws.Range("F17").Value = 36
$ws.Range("I17").Value = 191
$ws.Range("J17").Value = 201

# Row 19
$ws.Range("B19").Value = 8
$ws.Range("C19").Value = 10
$ws.Range("I19").Value = 8
$ws.Range("J19").Value = 10

# Row 20
$ws.Range("B20").Value = 183
$ws.Range("C20").Value = 191
$ws.Range("D20").Value = -0.041
$ws.Range("I20").Value = 183
$ws.Range("J20").Value = 191

# Row 22
$ws.Range("B22").Value = 36
$ws.Range("C22").Value = 36
$ws.Range("D22").Value = -0.021
$ws.Range("E22").Value = 36
$ws.Range("F22").Value = 36

# Row 23
$ws.Range("B23").Value = 577
$ws.Range("C23").Value = 593
$ws.Range("E23").Value = 390
$ws.Range("F23").Value = 379
$ws.Range("G23").Value = 169
$ws.Range("H23").Value = 192
$ws.Range("I23").Value = 18
$ws.Range("J23").Value = 22

# Row 26
$ws.Range("B26").Value = 577
$ws.Range("C26").Value = 593
$ws.Range("E26").Value = 390
$ws.Range("F26").Value = 379
$ws.Range("G26").Value = 169
$ws.Range("H26").Value = 192
$ws.Range("I26").Value = 18
$ws.Range("J26").Value = 22

# Row 31
$ws.Range("B31").Value = 4984
$ws.Range("C31").Value = 4952
$ws.Range("D31").Value = 0.007
$ws.Range("G31").Value = 4643
$ws.Range("H31").Value = 4572
$ws.Range("I31").Value = 341
$ws.Range("J31").Value = 380

# Row 34
$ws.Range("B34").Value = 3374
$ws.Range("C34").Value = 3299
$ws.Range("D34").Value = 0.023
$ws.Range("G34").Value = 3374
$ws.Range("H34").Value = 3299

# Row 36
$ws.Range("B36").Value = 732
$ws.Range("C36").Value = 738
$ws.Range("D36").Value = -0.009
$ws.Range("G36").Value = 732
$ws.Range("H36").Value = 738
$ws.Range("J36").Value = 0.2

# Row 39
$ws.Range("B39").Value = 878
$ws.Range("C39").Value = 914
$ws.Range("D39").Value = -0.039
$ws.Range("G39").Value = 537
$ws.Range("H39").Value = 534
$ws.Range("I39").Value = 341
$ws.Range("J39").Value = 380

# Row 46
$ws.Range("C46").Value = 7
$ws.Range("D46").Value = 0.03
$ws.Range("L46").Value = 7

# Row 49
$ws.Range("C49").Value = 7
$ws.Range("D49").Value = 0.03
$ws.Range("L49").Value = 7

# Row 51
$ws.Range("B51").Value = 2
$ws.Range("D51").Value = -0.17
$ws.Range("G51").Value = 2

# Row 58
$ws.Range("B58").Value = 2
$ws.Range("D58").Value = -0.17
$ws.Range("G58").Value = 2

# Row 60
$ws.Range("B60").Value = 588
$ws.Range("C60").Value = 718
$ws.Range("D60").Value = -0.18
$ws.Range("G60").Value = 588
$ws.Range("H60").Value = 718

# Row 61
$ws.Range("B61").Value = 384
$ws.Range("C61").Value = 472
$ws.Range("D61").Value = -0.19
$ws.Range("G61").Value = 384
$ws.Range("H61").Value = 472

# Row 62
$ws.Range("B62").Value = 81
$ws.Range("C62").Value = 106
$ws.Range("D62").Value = -0.23
$ws.Range("G62").Value = 81
$ws.Range("H62").Value = 106

# Row 63
$ws.Range("B63").Value = 123
$ws.Range("C63").Value = 141
$ws.Range("G63").Value = 123
$ws.Range("H63").Value = 141

# Row 64
$ws.Range("B64").Value = 408
$ws.Range("C64").Value = 350
$ws.Range("D64").Value = 0.17
$ws.Range("I64").Value = 408
$ws.Range("J64").Value = 350

# Row 66
$ws.Range("B66").Value = 408
$ws.Range("C66").Value = 350
$ws.Range("D66").Value = 0.17
$ws.Range("I66").Value = 408
$ws.Range("J66").Value = 350

# Row 67
$ws.Range("B67").Value = 14658
$ws.Range("C67").Value = 15160
$ws.Range("D67").Value = -0.033
$ws.Range("E67").Value = 425
$ws.Range("F67").Value = 415
$ws.Range("G67").Value = 12176
$ws.Range("H67").Value = 12565
$ws.Range("I67").Value = 2050
$ws.Range("J67").Value = 2173
$ws.Range("L67").Value = 7
